# Add the new "alcoholic" class specials (game_class_id = 10) to the
# "Class Specials" worksheet, appending rows 83-90 below the existing data
# (which currently ends at row 82).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 83 - Drinking Contest
$ws.Cells.Item(83, 1).Value = 82
$ws.Cells.Item(83, 2).Value = 10
$ws.Cells.Item(83, 3).Value = "Drinking Contest"
$ws.Cells.Item(83, 4).Value = "Lets have a drinking contest and see who gets blackout drunk first!"
$ws.Cells.Item(83, 5).Value = 1
$ws.Cells.Item(83, 9).Value = 0.05
$ws.Cells.Item(83, 13).Value = 0.1
$ws.Cells.Item(83, 14).Value = 0.05

# Row 84 - Alcoholics Rage
$ws.Cells.Item(84, 1).Value = 83
$ws.Cells.Item(84, 2).Value = 10
$ws.Cells.Item(84, 3).Value = "Alcoholics Rage"
$ws.Cells.Item(84, 4).Value = "You get black out drunk and flip tables, rage out and get extremely violent towards others!"
$ws.Cells.Item(84, 5).Value = 12
$ws.Cells.Item(84, 9).Value = 0.08
$ws.Cells.Item(84, 14).Value = 0.15

# Row 85 - Passed out drunk
$ws.Cells.Item(85, 1).Value = 84
$ws.Cells.Item(85, 2).Value = 10
$ws.Cells.Item(85, 3).Value = "Passed out drunk"
$ws.Cells.Item(85, 4).Value = "You are so drunk that you end up passing out. Nothing seems to be able to wake you. (Increases Armour and Health)"
$ws.Cells.Item(85, 5).Value = 24
$ws.Cells.Item(85, 10).Value = 0.1
$ws.Cells.Item(85, 13).Value = 0.2

# Row 86 - Blind Rage
$ws.Cells.Item(86, 1).Value = 85
$ws.Cells.Item(86, 2).Value = 10
$ws.Cells.Item(86, 3).Value = "Blind Rage"
$ws.Cells.Item(86, 4).Value = "Attack your foes with all you have! Drink your drink, rage on out! (Deals 1000 Damage, using 5% of your damage stat and growing by 10 over time for an additional 1000 damage)"
$ws.Cells.Item(86, 5).Value = 48
$ws.Cells.Item(86, 6).Value = 1000
$ws.Cells.Item(86, 7).Value = 10
$ws.Cells.Item(86, 8).Value = 0.05
$ws.Cells.Item(86, 9).Value = 0.2
$ws.Cells.Item(86, 13).Value = 0.3
$ws.Cells.Item(86, 14).Value = 0.25

# Row 87 - The bottom of the bottle
$ws.Cells.Item(87, 1).Value = 86
$ws.Cells.Item(87, 2).Value = 10
$ws.Cells.Item(87, 3).Value = "The bottom of the bottle"
$ws.Cells.Item(87, 4).Value = "Work your way to the bottom of the bottle of where your courage truly lies. Battle is upon you. (deal 10,000 Damage, using 10% of your damage stat and an additional 1000 damage for a total of an extra 10,000 damage on top.)"
$ws.Cells.Item(87, 5).Value = 60
$ws.Cells.Item(87, 6).Value = 10000
$ws.Cells.Item(87, 7).Value = 1000
$ws.Cells.Item(87, 8).Value = 0.1
$ws.Cells.Item(87, 9).Value = 0.2
$ws.Cells.Item(87, 10).Value = 0.1
$ws.Cells.Item(87, 13).Value = 0.5
$ws.Cells.Item(87, 14).Value = 0.4
$ws.Cells.Item(87, 15).Value = "attack"

# Row 88 - Fists of hate
$ws.Cells.Item(88, 1).Value = 87
$ws.Cells.Item(88, 2).Value = 10
$ws.Cells.Item(88, 3).Value = "Fists of hate"
$ws.Cells.Item(88, 4).Value = "Deal incredible damage with your fists as you pummel the enemy into a bloody mess. (Deal 40,000 in damage using 20% of your damage stat and growing by 4000 damage over time for an additional 40,000 damage)"
$ws.Cells.Item(88, 5).Value = 70
$ws.Cells.Item(88, 6).Value = 40000
$ws.Cells.Item(88, 7).Value = 4000
$ws.Cells.Item(88, 8).Value = 0.2
$ws.Cells.Item(88, 9).Value = 0.4
$ws.Cells.Item(88, 14).Value = 0.6
$ws.Cells.Item(88, 15).Value = "attack"

# Row 89 - The bottle is life
$ws.Cells.Item(89, 1).Value = 88
$ws.Cells.Item(89, 2).Value = 10
$ws.Cells.Item(89, 3).Value = "The bottle is life"
$ws.Cells.Item(89, 4).Value = "Alcohol gives you strength, gives you life. it's all you know. (Deal 50,000 damage using 30% of your damage stat and growing by 5,000 damage over time for an additional 50,000 damage)"
$ws.Cells.Item(89, 5).Value = 80
$ws.Cells.Item(89, 6).Value = 50000
$ws.Cells.Item(89, 7).Value = 5000
$ws.Cells.Item(89, 8).Value = 0.3
$ws.Cells.Item(89, 13).Value = 0.7
$ws.Cells.Item(89, 15).Value = "attack"

# Row 90 - Demons in your veins
$ws.Cells.Item(90, 1).Value = 89
$ws.Cells.Item(90, 2).Value = 10
$ws.Cells.Item(90, 3).Value = "Demons in your veins"
$ws.Cells.Item(90, 4).Value = "You drink and drink and cannot seem to quench the thirst, you cannot seem to quiet the demons in your head. (Deal 80,000 Damage, growing by 8,000 damage for an additional 80,000 damage and use 40% of your damage stat towards the damage)"
$ws.Cells.Item(90, 5).Value = 90
$ws.Cells.Item(90, 6).Value = 80000
$ws.Cells.Item(90, 7).Value = 8000
$ws.Cells.Item(90, 8).Value = 0.4
$ws.Cells.Item(90, 9).Value = 0.5
$ws.Cells.Item(90, 13).Value = 1
$ws.Cells.Item(90, 14).Value = 0.5
$ws.Cells.Item(90, 15).Value = "attack"
